$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 195.33333
$ws.Range("I9").Value = 115
$ws.Range("K9").Value = 115
$ws.Range("M9").Value = 54
# Row 53
$ws.Range("H53").Value = 204.72223
$ws.Range("I53").Value = 59.214287
$ws.Range("K53").Value = 59.214287
$ws.Range("M53").Value = 577.785713
# Row 62
$ws.Range("H62").Value = 5834.5
$ws.Range("I62").Value = 2680.6
$ws.Range("J62").Value = 7268.091
$ws.Range("K62").Value = 2680.6
$ws.Range("L62").Value = 7268.091
$ws.Range("M62").Value = -2056.6
$ws.Range("N62").Value = -8516.091
# Row 65
$ws.Range("H65").Value = 5834.5
$ws.Range("I65").Value = 2680.6
$ws.Range("J65").Value = 7268.091
$ws.Range("K65").Value = 13403
$ws.Range("L65").Value = 36340.455
$ws.Range("M65").Value = -10283
$ws.Range("N65").Value = -42580.455
# Row 74
$ws.Range("H74").Value = 4465.625
$ws.Range("I74").Value = 3318.182
$ws.Range("K74").Value = 3318.182
$ws.Range("M74").Value = -2382.182
# Row 77
$ws.Range("H77").Value = 4465.625
$ws.Range("I77").Value = 3318.182
$ws.Range("K77").Value = 16590.91
$ws.Range("M77").Value = -11910.91
# Row 97
$ws.Range("H97").Value = 2095.2727
$ws.Range("J97").Value = 1462.8
$ws.Range("L97").Value = 4388.4
$ws.Range("N97").Value = -5380.4
# Row 98
$ws.Range("H98").Value = 838
$ws.Range("I98").Value = 838
$ws.Range("K98").Value = 838
$ws.Range("M98").Value = 660
# Row 112
$ws.Range("H112").Value = 3939.1538
$ws.Range("J112").Value = 3860.9
$ws.Range("L112").Value = 11582.7
$ws.Range("N112").Value = -13798.7
# Row 122
$ws.Range("H122").Value = 838
$ws.Range("I122").Value = 838
$ws.Range("K122").Value = 2514
$ws.Range("M122").Value = -64
# Row 125
$ws.Range("H125").Value = 8426.700000000001
$ws.Range("I125").Value = 7782.6665
$ws.Range("J125").Value = 9392.75
$ws.Range("K125").Value = 70043.9985
$ws.Range("L125").Value = 84534.75
$ws.Range("M125").Value = -67583.9985
$ws.Range("N125").Value = -89454.75
# Row 137
$ws.Range("H137").Value = 8898.875
$ws.Range("I137").Value = 9813
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 29439
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = -26889
$ws.Range("N137").Value = -12600

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3117.5
$ws.Range("I45").Value = 2236
$ws.Range("K45").Value = 2236
$ws.Range("M45").Value = -1859
# Row 97
$ws.Range("H97").Value = 741.8889
$ws.Range("I97").Value = 489.7857
$ws.Range("J97").Value = 1624.25
$ws.Range("K97").Value = 489.7857
$ws.Range("L97").Value = 1624.25
$ws.Range("M97").Value = 6.21429999999998
$ws.Range("N97").Value = -2616.25
# Row 101
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 986.7857
$ws.Range("I7").Value = 983.1818
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 983.1818
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = -870.1818
$ws.Range("N7").Value = -1226
# Row 21
$ws.Range("H21").Value = 1938.3334
# Row 31
$ws.Range("H31").Value = 12157.833
$ws.Range("I31").Value = 17298.5
$ws.Range("K31").Value = 17298.5
$ws.Range("M31").Value = -17003.5
# Row 34
$ws.Range("H34").Value = 12157.833
$ws.Range("I34").Value = 17298.5
$ws.Range("K34").Value = 17298.5
$ws.Range("M34").Value = -17096.5
# Row 74
$ws.Range("H74").Value = 32733.166
$ws.Range("J74").Value = 32733.166
$ws.Range("L74").Value = 32733.166
$ws.Range("N74").Value = -34481.166
# Row 77
$ws.Range("H77").Value = 32733.166
$ws.Range("J77").Value = 32733.166
$ws.Range("L77").Value = 98199.49800000001
$ws.Range("N77").Value = -106935.498
# Row 106
$ws.Range("H106").Value = 25994.5
$ws.Range("J106").Value = 25994.5
$ws.Range("L106").Value = 25994.5
$ws.Range("N106").Value = -28518.5

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 55
$ws.Range("H55").Value = 5066.6665
$ws.Range("J55").Value = 5225
$ws.Range("L55").Value = 15675
$ws.Range("N55").Value = -16029

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 18
$ws.Range("H18").Value = 4102.5
$ws.Range("I18").Value = 8005
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 8005
$ws.Range("L18").Value = 200
$ws.Range("M18").Value = -7712
$ws.Range("N18").Value = -786
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
# Row 31
$ws.Range("H31").Value = 3366.2
$ws.Range("I31").Value = 2957.75
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 2957.75
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -2665.75
$ws.Range("N31").Value = -5584
# Row 37
$ws.Range("H37").Value = 3366.2
$ws.Range("I37").Value = 2957.75
$ws.Range("J37").Value = 5000
$ws.Range("K37").Value = 2957.75
$ws.Range("L37").Value = 5000
$ws.Range("M37").Value = -2680.75
$ws.Range("N37").Value = -5554
# Row 50
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
# Row 80
$ws.Range("H80").Value = 3038.889
$ws.Range("I80").Value = 2982.353
$ws.Range("K80").Value = 2982.353
$ws.Range("M80").Value = -1984.353
# Row 83
$ws.Range("H83").Value = 3038.889
$ws.Range("I83").Value = 2982.353
$ws.Range("K83").Value = 14911.765
$ws.Range("M83").Value = -9919.764999999999

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 4
$ws.Range("H4").Value = 6750.75
$ws.Range("I4").Value = 3627
$ws.Range("J4").Value = 9874.5
$ws.Range("K4").Value = 3627
$ws.Range("L4").Value = 9874.5
$ws.Range("M4").Value = -3514
$ws.Range("N4").Value = -10100.5
# Row 19
$ws.Range("H19").Value = 4500
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
# Row 25
$ws.Range("H25").Value = 8834.223
$ws.Range("J25").Value = 8334.666999999999
$ws.Range("L25").Value = 8334.666999999999
$ws.Range("N25").Value = -8794.666999999999
# Row 28
$ws.Range("H28").Value = 6750.75
$ws.Range("I28").Value = 3627
$ws.Range("J28").Value = 9874.5
$ws.Range("K28").Value = 3627
$ws.Range("L28").Value = 9874.5
$ws.Range("M28").Value = -3395
$ws.Range("N28").Value = -10338.5
# Row 37
$ws.Range("H37").Value = 6750.75
$ws.Range("I37").Value = 3627
$ws.Range("J37").Value = 9874.5
$ws.Range("K37").Value = 3627
$ws.Range("L37").Value = 9874.5
$ws.Range("M37").Value = -3520
$ws.Range("N37").Value = -10088.5
# Row 46
$ws.Range("H46").Value = 3500
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 5000
$ws.Range("N46").Value = -5376
# Row 136
$ws.Range("H136").Value = 7534.923
$ws.Range("I136").Value = 6518.4443
$ws.Range("K136").Value = 19555.3329
$ws.Range("M136").Value = -17005.3329

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1219
$ws.Range("I126").Value = 1268.8182
$ws.Range("J126").Value = 945
$ws.Range("K126").Value = 3806.4546
$ws.Range("L126").Value = 2835
$ws.Range("M126").Value = -1336.4546
$ws.Range("N126").Value = -7775
